$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reporting-period dates in row 8 (Q3 2021 -> Q4 2021 refresh)
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561
$ws.Range("E8").Value = 44571
$ws.Range("H8").Value = 44571
$ws.Range("I8").Value = 44571

# Move the view/selection to match where the editor left off
$ws.Activate() | Out-Null
$ws.Range("J8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 6
